# Add a "Skill Description" column (full skill/attribute name) as the new
# column B, shifting the existing SFIA Level / Keycode / Description columns
# one place to the right (C/D/E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, empty column at B - this shifts old B->C, C->D, D->E.
$ws.Columns("B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Skill Description"

# Map each SkillCode / attribute code (column A) to its full descriptive name
# that now belongs in column B.
$fullNames = @{
    "Autonomy"   = "Autonomy";
    "Influence"  = "Influence";
    "Complexity" = "Complexity";
    "Knowledge"  = "Knowledge";
    "POMG"       = "Portfolio management";
    "PGMG"       = "Programme management";
    "PRMG"       = "Project management";
    "RLMT"       = "Stakeholder relationship management";
    "ISCO"       = "Information systems coordination";
}

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $code = $ws.Cells.Item($row, 1).Value2
    if ($code -and $fullNames.ContainsKey($code)) {
        $ws.Cells.Item($row, 2).Value = $fullNames[$code]
    }
}
